$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing G1 (value "B", shared string) moves to J1.
# New cells inserted: E1=5, F1=2 (old E1=2 shifts to F1), G1=7 (new), H1=3 (old F1 shifts to H1), I1=2 (new)
$ws.Range("J1").Value = $ws.Range("G1").Value2
$ws.Range("H1").Value = 3
$ws.Range("E1").Value = 5
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 7
$ws.Range("I1").Value = 2

$ws.Range("E1:J1").HorizontalAlignment = -4108

$ws.Range("H20").Select()
